$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51 (columns B=Coin, C=Link, D=Price, E=Volume(1h))
$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "29.780.00", "  +2.07%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.856.93", "  +1.65%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9997", "  +0.03%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "245.18", "  +1.35%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.6407", "  +3.44%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  -0.03%  "),
    @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "47.09", "  +4.73%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07499", "  +2.18%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2974", "  +2.66%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "24.19", "  +4.95%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07684", "  +0.01%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.875.63", "  +2.70%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.059", "  +1.99%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6852", "  +3.07%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "83.98", "  +2.01%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000009501", "  +6.02%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.085", "  +3.70%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "29.763.10", "  +2.12%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.120.89", "  +2.40%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "239.70", "  +0.74%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "12.68", "  +1.91%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9998", "  -0.06%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.427", "  +1.93%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  +0.05%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "158.90", "  +0.39%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1433", "  +0.40%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.522", "  +0.39%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "17.95", "  +1.54%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06186", "  +9.76%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.502", "  +1.52%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.276", "  +5.33%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.146", "  +1.20%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.108", "  +0.43%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.877", "  +2.18%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.160", "  +2.48%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7317", "  -0.44%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.604", "  -0.84%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.858", "  +0.52%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01796", "  +1.57%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.212.72", "  -0.23%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9245", "  +0.43%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.184", "  -2.04%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.001", "  -0.02%  "),
    @("RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "2.025.12", "  +2.66%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "102.02", "  +0.22%  "),
    @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "66.36", "  +2.24%  "),
    @("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.00000000121", "  +4.15%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.349", "  +2.60%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4082", "  +1.19%  "),
    @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1132", "  +1.47%  ")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row++
}
